# Auto-save inventory data - 2025-07-04 06:48:34
#
# Adds a Model/SN-Lot entry to the last inventory row (row 23):
#   D23 = "Test"   (Model)
#   E23 = "123"    (SN/Lot)
#
# "123" must be written as literal text (not the number 123), matching the
# source OOXML diff where it is stored as a new shared string. A direct
# Range.Value = "123" assignment would be auto-coerced to a number by the
# COM layer's numeric-literal detection, so we build the text via a
# formula that yields a string ("="123"""), then copy/paste-special just
# the resulting value back into the target cell. That keeps the cell's
# type as text without touching any cell styles/number formats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D23").Value = "Test"

$helper = $ws.Range("Z1")
$helper.Formula = "=""123"""
$helper.Copy()
$ws.Range("E23").PasteSpecial(-4163)
$helper.ClearContents()
